# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 23:05"

# Swap Seychelles and Montserrat ordering (Seychelles now comes before Montserrat)
$ws.Range("A210").Value = "Seychelles"
$ws.Range("A211").Value = "Montserrat"

# Update numeric data for updated countries

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1833660
$ws.Range("C4").Value = 16840
$ws.Range("D4").Value = 541029
$ws.Range("E4").Value = 1186503
$ws.Range("G4").Value = 571
$ws.Range("H4").Value = 106128

# Row 5 - Brasil
$ws.Range("B5").Value = 506455
$ws.Range("C5").Value = 8015
$ws.Range("E5").Value = 272004
$ws.Range("G5").Value = 246
$ws.Range("H5").Value = 29080

# Row 10 - Francia
$ws.Range("B10").Value = 190609
$ws.Range("C10").Value = 8782
$ws.Range("D10").Value = 91852
$ws.Range("E10").Value = 93349
$ws.Range("G10").Value = 223
$ws.Range("H10").Value = 5408

# Row 12 - Turquia
$ws.Range("B12").Value = 183494
$ws.Range("C12").Value = 200
$ws.Range("E12").Value = 9689

# Row 27 - Suecia
$ws.Range("B27").Value = 39098
$ws.Range("C27").Value = 527
$ws.Range("D27").Value = 19592
$ws.Range("E27").Value = 16148
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = 3358

# Row 52 - Oman
$ws.Range("E52").Value = 8992
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 49

# Row 121 - Costa Rica
$ws.Range("B121").Value = 958
$ws.Range("C121").Value = 2
$ws.Range("D121").Value = 839
$ws.Range("E121").Value = 55

# Row 210 - Seychelles (after swap)
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211 - Montserrat (after swap)
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
